# merge the mainicon and sceneconfig
# Removes the obsolete "休息" (Rest) MainIcon row (row 19) from the sheet;
# all subsequent rows shift up by one and the table/dimension shrink accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(19).Delete()

$ws.Range("E20").Select()
